$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new columns (B:J) before the existing data, shifting old B:V to K:AE.
$ws.Columns("B:J").Insert()

# Populate the new header row cells with the newer weekly snapshot labels.
# Written right-to-left so the shared-string table grows in chronological
# order (Jun_16 first ... Sep_08 last), matching the source workbook.
$ws.Range("J1").Value = "Jun_16"
$ws.Range("I1").Value = "Jun_24"
$ws.Range("H1").Value = "Jun_30"
$ws.Range("G1").Value = "Jul_07"
$ws.Range("F1").Value = "Jul_17"
$ws.Range("E1").Value = "Jul_23"
$ws.Range("D1").Value = "Aug_04"
$ws.Range("C1").Value = "Aug_25"
$ws.Range("B1").Value = "Sep_08"

# The new snapshot columns carry the same "UN" (unchanged) rating used
# throughout the rest of the table for every data row.
$ws.Range("B2:J33").Value = "UN"

# Restore the column widths for the newly inserted columns (same width as
# the rest of the data columns) and make sure the formatting of the
# trailing columns matches as well.
$ws.Range("B1:AE1").EntireColumn.ColumnWidth = 46.1640625
